$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116, shifting existing row 116 (and below) down to 117.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new data record.
$ws.Range("A116").Value = 10
$ws.Range("B116").Value = "Vega Modelo de Temuco"
$ws.Range("C116").Value = "La Araucanía"
$ws.Range("D116").Value = "2022-05-24"
$ws.Range("E116").Value = 9
$ws.Range("F116").Value = 100112052
$ws.Range("G116").Value = "Albahaca"
$ws.Range("H116").Value = "Sin especificar"
$ws.Range("I116").Value = "Primera"
$ws.Range("J116").Value = 30
$ws.Range("K116").Value = 5000
$ws.Range("L116").Value = 5000
$ws.Range("M116").Value = 5000
$ws.Range("N116").Value = "$/paquete"
$ws.Range("O116").Value = "Región de Arica y Parinacota"
$ws.Range("P116").Value = 5000
$ws.Range("Q116").Value = 1
$ws.Range("R116").Value = "Hortaliza"
